# Add 2022-Q3 data: the existing "2022-Q1" sheet becomes "2022-Q3" (new
# fund data), and a brand-new "2022-Q1" sheet is appended right after it
# holding the data that used to live in the "2022-Q1" sheet. The "总计"
# (total) sheet gets a new row on top for the 2022-Q3 totals, pushing the
# old 2022-Q1 total row down.

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item(1)
$wsQ = $wb.Worksheets.Item(2)

# --- 1. Duplicate the current "2022-Q1" sheet so its fund-holding data is
#        preserved verbatim (values + formatting) under the "2022-Q1" name,
#        placed right after the sheet that will become "2022-Q3".
$wsQ.Copy($null, $wsQ)
$wsQ1New = $wb.Worksheets.Item(3)

# --- 2. Turn the original sheet into "2022-Q3" (rename it out of the way
#        first so the freshly-copied sheet can take the "2022-Q1" name) and
#        overwrite its data with the new quarter's fund row.
$wsQ.Name = "2022-Q3"
$wsQ1New.Name = "2022-Q1"

$wsQ.Range("B2:G2").NumberFormat = "@"

$wsQ.Range("B2").Value = "010714"
$wsQ.Range("C2").Value = "东方红远见价值混合"
$wsQ.Range("D2").Value = "15.24"
$wsQ.Range("E2").Value = "94.15"
$wsQ.Range("F2").Value = "4.36"
$wsQ.Range("G2").Value = "0.6645"
$wsQ.Range("H2").Value = 6

$wsQ.Range("B2:G2").Style = "Normal"

# --- 3. Insert the new 2022-Q3 summary row at the top of the "总计" sheet,
#        pushing the existing 2022-Q1 row down to row 3.
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Rows.Item(2).ClearFormats()

# restore the row-2 "index" column formatting (it matches row 3's, which
# kept its original style when it was shifted down).
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 0.66

$wsTotal.Range("A3").Value = 1
